# Apply updated values to the "RF" worksheet of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RF")

# Row 13 - Gross Margin
$ws.Range("D13").Value = 1.0169
$ws.Range("E13").Value = 0.9674
$ws.Range("F13").Value = 0.9154
$ws.Range("G13").Value = 0.8732

# Row 15 - EBT margin
$ws.Range("D15").Value = 0.1779
$ws.Range("E15").Value = 0.1542
$ws.Range("F15").Value = 0.2656
$ws.Range("G15").Value = 0.2957

# Row 16 - Net Profit Margin
$ws.Range("D16").Value = 0.1286
$ws.Range("E16").Value = 0.1065
$ws.Range("F16").Value = 0.1987
$ws.Range("G16").Value = 0.2239

# Row 17 - Free Cash Flow Margin
$ws.Range("D17").Value = 0.5002
$ws.Range("E17").Value = 0.4996
$ws.Range("F17").Value = 0.3897
$ws.Range("G17").Value = 0.3544

# Row 27 - Operating Cash Flow Margin
$ws.Range("D27").Value = 0.5314
$ws.Range("E27").Value = 0.5389
$ws.Range("F27").Value = 0.4219
$ws.Range("G27").Value = 0.3845
